$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells whose new value looks like a plain number need to be
# forced to Text format first, otherwise Excel silently converts the
# assigned string into a floating point number (losing exact formatting,
# e.g. "23.30" -> 23.3, or "0.0270" -> 2.7E-2).

$ws.Range("D2").Value = "40.182.25"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "2.229.91"
$ws.Range("E3").Value = "  +0.77%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.73"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  -0.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.72"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.95"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +6.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0783"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("E13").Value = "  +2.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.46"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").Value = "2.576.65"
$ws.Range("E15").Value = "  +0.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.84"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.19%  "

$ws.Range("D17").Value = "2.221.16"
$ws.Range("E17").Value = "  +0.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.737"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.19%  "

$ws.Range("D19").Value = "40.117.22"
$ws.Range("E19").Value = "  +0.42%  "

$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.26"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.79"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.81"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.75"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.36%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.48"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.83"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.30"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.17%  "

$ws.Range("E29").Value = "  -2.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.32"
$ws.Range("D30").ClearFormats()

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.51"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.79"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.23%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.98"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.06"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0717"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.31"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.43%  "

$ws.Range("E38").Value = "  +1.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.76"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0995"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.50"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.50%  "

$ws.Range("D42").Value = "2.095.10"
$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.71"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.79"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +6.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.13"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0270"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.35%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.96"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -10.99%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.72"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.36%  "

$ws.Range("D49").Value = "2.450.76"
$ws.Range("E49").Value = "  +0.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.48"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.51%  "

$ws.Range("E51").Value = "  +3.53%  "
